# Hortaliza, Femacal de La Calera - Alcachofa
# Weekly update: insert two new price-report rows (new rows 227 and 228),
# pushing all existing data (previously rows 227-320) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 227.
$ws.Rows("227:228").Insert()

# --- New row 227 ---
$ws.Range("A227").Value = 3
$ws.Range("B227").Value = "Femacal de La Calera"
$ws.Range("C227").Value = "Coquimbo"
$ws.Range("D227").Value = 44726
$ws.Range("E227").Value = 5
$ws.Range("F227").Value = 100112013
$ws.Range("G227").Value = "Alcachofa"
$ws.Range("H227").Value = "Argentina(o)"
$ws.Range("I227").Value = "Primera"
$ws.Range("J227").Value = 85
$ws.Range("K227").Value = 16000
$ws.Range("L227").Value = 16500
$ws.Range("M227").Value = 16265
$ws.Range("N227").Value = "`$/caja 50 unidades"
$ws.Range("O227").Value = "Provincia de Limarí"
$ws.Range("P227").Value = 325
$ws.Range("Q227").Value = 50
$ws.Range("R227").Value = "Hortaliza"

# --- New row 228 ---
$ws.Range("A228").Value = 3
$ws.Range("B228").Value = "Femacal de La Calera"
$ws.Range("C228").Value = "Coquimbo"
$ws.Range("D228").Value = 44726
$ws.Range("E228").Value = 5
$ws.Range("F228").Value = 100112013
$ws.Range("G228").Value = "Alcachofa"
$ws.Range("H228").Value = "Española"
$ws.Range("I228").Value = "Primera"
$ws.Range("J228").Value = 42
$ws.Range("K228").Value = 21000
$ws.Range("L228").Value = 21000
$ws.Range("M228").Value = 21000
$ws.Range("N228").Value = "`$/caja 50 unidades"
$ws.Range("O228").Value = "Provincia de Limarí"
$ws.Range("P228").Value = 420
$ws.Range("Q228").Value = 50
$ws.Range("R228").Value = "Hortaliza"
